# Update metadata and running round4 imaging
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# --- 1. Expand Table1 to add a new "ShortName" column, inserted just before "Comment".
#     The engine's ListColumns.Add() always appends, so we grow the table by one
#     column (old "Comment" @ N, new blank @ O), then shift the Comment data from
#     N into O and retarget N as "ShortName".
$tbl.Resize($ws.Range("A1:O30"))

for ($r = 1; $r -le 30; $r++) {
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 14).Value2
}

$ws.Cells.Item(1, 14).Value = "ShortName"
$ws.Cells.Item(1, 15).Value = "Comment"

# Clear out every old "Comment" cell in N (now ShortName) - we'll repopulate below.
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 14).ClearContents()
}

# Only rows 6 (Test1) and 27 (CM3) keep a Comment; clear the rest that got copied over.
for ($r = 2; $r -le 30; $r++) {
    if ($r -ne 6 -and $r -ne 27) {
        $ws.Cells.Item($r, 15).ClearContents()
    }
}
$ws.Cells.Item(27, 15).Value = "Problem with Pulse, spots"

# --- 2. ShortName values (new column N) for every data row.
$shortNames = @{
    2  = "Z552_1"
    3  = "Z552_2"
    4  = "Z673_1"
    5  = "Z673_2"
    6  = "Test1"
    7  = "Z673_3"
    8  = "Z552_3"
    9  = "Z552_4"
    10 = "BM24"
    11 = "Z552_5"
    12 = "Z673_4"
    13 = "BM26"
    14 = "BM27"
    15 = "Z673_5"
    16 = "BM28"
    17 = "MC2_A"
    18 = "MC2_B"
    19 = "MC2_C"
    20 = "MC2_D"
    21 = "MC2_E"
    22 = "MC2_F"
    23 = "BM29"
    24 = "BM30"
    25 = "CM1"
    26 = "CM2"
    27 = "CM3"
    28 = "CF1"
    29 = "CF2"
    30 = "CF3"
}
foreach ($r in $shortNames.Keys) {
    $ws.Cells.Item($r, 14).Value = $shortNames[$r]
}

# --- 3. Imaging round 4 is underway: the placeholder "Not imaged yet" / "Not
#     perfused yet" Path entries (rows 23-30) now point at the round-4 folder.
$newRoundPath = "/nrs/spruston/Boaz/I2/20241104_iDISCO_R4/"
for ($r = 23; $r -le 30; $r++) {
    $ws.Cells.Item($r, 10).Value = $newRoundPath
}

# --- 4. Column widths / formatting touch-ups to roughly match the new layout.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(12).ColumnWidth
$ws.Columns.Item(15).ColumnWidth = 28.33

# --- 5. Selection moved to N16 and the frozen/top-left scroll reset to the top.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("N16").Select()
